# Automatische test-sync: 2025-08-03 15:04:50
# Adds a new test-mail log row (#15) to the "Logs" sheet and refreshes the
# dependent "Planning / Afspraak" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append the new log entry as row 23 -----------------------------------
$row = 23
$logs.Cells.Item($row, 1).Value = "Leg dit even neer bij Koen."
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #15: Leg dit even neer bij Koen."
$logs.Cells.Item($row, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($row, 6).Value = "2025-08-03 15:04:15"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

# --- Extend the conditional-formatting ranges so row 23 is included -------
$logs.Range("D2:D22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D23"))
$logs.Range("G2:G22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G23"))
$logs.Range("H2:H22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H23"))
$logs.Range("I2:I22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I23"))
$logs.Range("J2:J22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J23"))

# --- Refresh the Dashboard summary count for "Planning / Afspraak" --------
$dash.Range("B4").Value = 5
